# "Moly and Pax - Weekly => Monthly"
# Rename the Aged Care Molnupiravir / Paxlovid prescription metrics from
# "(Weekly)" to "(Monthly)" on the Metrics sheet, update the active
# selection, and refresh the page setup (A4 portrait) to match the
# author's re-saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metrics")

# Rows 72-77 hold the six Molnupiravir/Paxlovid prescription metric labels
# (column C). Update each in place so the underlying shared-string slot is
# reused rather than a new one being appended.
$ws.Range("C72").Value = "# Aged Care Molnupiravir Prescriptions (Monthly)"
$ws.Range("C73").Value = "# Aged Care Molnupiravir Prescriptions (Monthly) per 1M"
$ws.Range("C74").Value = "% Aged Care Molnupiravir Prescriptions (Monthly) per Case"
$ws.Range("C75").Value = "# Aged Care Paxlovid Prescriptions (Monthly)"
$ws.Range("C76").Value = "# Aged Care Paxlovid Prescriptions (Monthly) per 1M"
$ws.Range("C77").Value = "% Aged Care Paxlovid Prescriptions (Monthly) per Case"

# Match the page setup captured in the re-saved file (A4 portrait).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# The author's cursor ended up on C73 when the file was saved.
$ws.Range("C73").Select()

Write-Output "Updated Molnupiravir/Paxlovid (Weekly -> Monthly) labels"
